$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 27, duplicating row 26 (Case File - creator read access)
# but for Complaint.
$ws.Range("A26:H26").Copy()
$ws.Range("A27:H27").PasteSpecial()

$ws.Cells.Item(27, 2).Value = "Complaint - creator read access"
$ws.Cells.Item(27, 3).Value = "COMPLAINT"
$ws.Cells.Item(27, 8).Value = "reader, creator"

$ws.Range("B28").Select()
